$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the list of D-column cells whose new values look numeric; Excel
# would otherwise silently convert those assignments from Text to Number,
# which would change the cell type from the original inlineStr/text cells.
# Force those ranges to Text format first, assign, then restore the style
# so no residual formatting diff is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '62.273.12'
$ws.Range("D3").Value = '2.453.67'
$ws.Range("D5").Value = '580.87'
$ws.Range("D6").Value = '143.74'
$ws.Range("D8").Value = '0.531'
$ws.Range("D9").Value = '2.452.73'
$ws.Range("D10").Value = '0.107'
$ws.Range("D12").Value = '5.21'
$ws.Range("D13").Value = '0.345'
$ws.Range("D14").Value = '26.54'
$ws.Range("D15").Value = '0.0000173'
$ws.Range("D16").Value = '2.824.13'
$ws.Range("D17").Value = '62.176.01'
$ws.Range("D18").Value = '2.441.37'
$ws.Range("D19").Value = '10.90'
$ws.Range("D20").Value = '7.17'
$ws.Range("D21").Value = '329.80'
$ws.Range("D22").Value = '4.10'
$ws.Range("D23").Value = '2.00'
$ws.Range("D25").Value = '66.03'
$ws.Range("D26").Value = '9.42'
$ws.Range("D27").Value = '620.81'
$ws.Range("D28").Value = '0.0₃0956'
$ws.Range("D29").Value = '2.541.35'
$ws.Range("D31").Value = '1.43'
$ws.Range("D32").Value = '8.03'
$ws.Range("D34").Value = '1.88'
$ws.Range("D35").Value = '4.91'
$ws.Range("D38").Value = '0.377'
$ws.Range("D39").Value = '5.33'
$ws.Range("D40").Value = '149.54'
$ws.Range("D41").Value = '18.38'
$ws.Range("D45").Value = '2.46'
$ws.Range("D46").Value = '143.47'
$ws.Range("D47").Value = '3.65'
$ws.Range("D48").Value = '0.0525'
$ws.Range("D49").Value = '0.605'
$ws.Range("D50").Value = '19.59'
$ws.Range("D51").Value = '0.0₆0240'

$dRange.Style = "Normal"

$ws.Range("E2").Value = '  -1.54%  '
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -1.24%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -3.99%  '
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("E12").Value = '  -1.18%  '
$ws.Range("E13").Value = '  -3.21%  '
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("E15").Value = '  -4.10%  '
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("E19").Value = '  -3.78%  '
$ws.Range("E20").Value = '  -2.12%  '
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("E23").Value = '  -3.88%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("E26").Value = '  +6.07%  '
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("E28").Value = '  -6.84%  '
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  -5.15%  '
$ws.Range("E32").Value = '  -2.48%  '
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("E35").Value = '  -5.61%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -6.63%  '
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("E40").Value = '  +1.63%  '
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("E42").Value = '  -2.80%  '
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("E45").Value = '  -5.75%  '
$ws.Range("E46").Value = '  -3.69%  '
$ws.Range("E47").Value = '  -3.48%  '
$ws.Range("E48").Value = '  -1.57%  '
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("E50").Value = '  -7.68%  '
$ws.Range("E51").Value = '  +10.12%  '
